$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row with the country "jona" right after the existing
# CountryName / abc / def / ghi rows (A1:A4 -> A1:A5).
$cell = $ws.Range("A5")
$cell.Value = "jona"

# Select the newly added cell, matching the saved sheetView selection.
$cell.Select()

# The new row's cell carries a distinct cell style (a second <xf> entry
# was added to cellXfs in the saved workbook), so give it its own format.
$cell.Font.Bold = $true
